$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 14 ("8. Indica que a adição foi feita com sucesso"), which shifts
# everything below it up by one row (rows 15,16 -> 14,15) and shrinks the
# B6:B14 / B15:B16 merges accordingly.
$ws.Rows("14:14").Delete()

# Update the post-condition text.
$ws.Range("C5").Value2 = "Stock de uma peça foi alterada"

# Row 12: "6. Adiciona peça a carros que precisem da mesma" -> merged step text.
$ws.Range("D12").Value2 = "6. Adiciona Stock e adiciona peça a carros que precisem da mesma"

# Row 13: "7. Regista adição do stock" -> "7. Indica que a adição  foi feita com sucesso"
$ws.Range("D13").Value2 = "7. Indica que a adição  foi feita com sucesso"

# Row 15 (previously row 16 before the deletion): update exception message text.
$ws.Range("D15").Value2 = "4.1. Apresenta mensagem ""Peça Inválida"""

# Widen column D (target stored width is 64.875 character-units; the
# runtime quantizes ColumnWidth writes to steps of 1/6, so 64 is the closest
# achievable value, yielding a stored width of 64.83333333333333).
$ws.Columns("D").ColumnWidth = 64

# Update the saved selection to C5:D5.
$ws.Range("C5:D5").Select()
